$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 485.04166
$ws.Range("I28").Value = 346.6111
$ws.Range("J28").Value = 900.3333
$ws.Range("K28").Value = 346.6111
$ws.Range("L28").Value = 900.3333
$ws.Range("M28").Value = 138.3889
$ws.Range("N28").Value = -1870.3333

$ws.Range("H32").Value = 1517.2
$ws.Range("I32").Value = 1200
$ws.Range("J32").Value = 1834.4
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 1834.4
$ws.Range("M32").Value = -874
$ws.Range("N32").Value = -2486.4

$ws.Range("H111").Value = 125001270
$ws.Range("I111").Value = 250001120
$ws.Range("J111").Value = 1416
$ws.Range("K111").Value = 750003360
$ws.Range("L111").Value = 4248
$ws.Range("M111").Value = -750000293
$ws.Range("N111").Value = -10382

$ws.Range("H116").Value = 7695888
$ws.Range("I116").Value = 7695888
$ws.Range("K116").Value = 7695888
$ws.Range("M116").Value = -7692446

$ws.Range("H129").Value = 974.38666
$ws.Range("I129").Value = 475
$ws.Range("J129").Value = 988.0685
$ws.Range("K129").Value = 1425
$ws.Range("L129").Value = 2964.2055
$ws.Range("M129").Value = 3575
$ws.Range("N129").Value = -12964.2055

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3528.25
$ws.Range("I2").Value = 6500
$ws.Range("J2").Value = 2537.6667
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 2537.6667
$ws.Range("M2").Value = -6387
$ws.Range("N2").Value = -2763.6667

$ws.Range("H45").Value = 2338.3157
$ws.Range("I45").Value = 1799.091
$ws.Range("J45").Value = 3079.75
$ws.Range("K45").Value = 1799.091
$ws.Range("L45").Value = 3079.75
$ws.Range("M45").Value = -1422.091
$ws.Range("N45").Value = -3833.75

$ws.Range("H116").Value = 3528.25
$ws.Range("I116").Value = 6500
$ws.Range("J116").Value = 2537.6667
$ws.Range("K116").Value = 6500
$ws.Range("L116").Value = 2537.6667
$ws.Range("M116").Value = -4206
$ws.Range("N116").Value = -7125.6667

$ws.Range("H132").Value = 2675.5293
$ws.Range("I132").Value = 2528.6365
$ws.Range("J132").Value = 2944.8333
$ws.Range("K132").Value = 7585.9095
$ws.Range("L132").Value = 8834.499899999999
$ws.Range("M132").Value = -5055.9095
$ws.Range("N132").Value = -13894.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3528.25
$ws.Range("I3").Value = 6500
$ws.Range("J3").Value = 2537.6667
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 2537.6667
$ws.Range("M3").Value = -6386
$ws.Range("N3").Value = -2765.6667

$ws.Range("H80").Value = 108
$ws.Range("J80").Value = 79.59999999999999
$ws.Range("L80").Value = 79.59999999999999
$ws.Range("N80").Value = -2075.6

$ws.Range("H83").Value = 108
$ws.Range("J83").Value = 79.59999999999999
$ws.Range("L83").Value = 398
$ws.Range("N83").Value = -10382

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6744.4736
$ws.Range("I31").Value = 7023.1
$ws.Range("K31").Value = 7023.1
$ws.Range("M31").Value = -6728.1

$ws.Range("H34").Value = 6744.4736
$ws.Range("I34").Value = 7023.1
$ws.Range("K34").Value = 7023.1
$ws.Range("M34").Value = -6821.1

$ws.Range("H99").Value = 2988.9473
$ws.Range("I99").Value = 2585
$ws.Range("J99").Value = 4120
$ws.Range("K99").Value = 2585
$ws.Range("L99").Value = 4120
$ws.Range("M99").Value = -1087
$ws.Range("N99").Value = -7116

$ws.Range("H126").Value = 2988.9473
$ws.Range("I126").Value = 2585
$ws.Range("J126").Value = 4120
$ws.Range("K126").Value = 7755
$ws.Range("L126").Value = 12360
$ws.Range("M126").Value = -5285
$ws.Range("N126").Value = -17300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 900
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1889
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 900
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4044
$ws.Range("N71").ClearContents()

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H92").Value = 1428838.9
$ws.Range("I92").Value = 2500175.5
$ws.Range("J92").Value = 390
$ws.Range("K92").Value = 7500526.5
$ws.Range("L92").Value = 1170
$ws.Range("M92").Value = -7499278.5
$ws.Range("N92").Value = -3666

$ws.Range("H97").Value = 1080.2727
$ws.Range("I97").Value = 397.66666
$ws.Range("J97").Value = 1336.25
$ws.Range("K97").Value = 1192.99998
$ws.Range("L97").Value = 4008.75
$ws.Range("M97").Value = -696.9999800000001
$ws.Range("N97").Value = -5000.75

$ws.Range("H98").Value = 215.7
$ws.Range("I98").Value = 217
$ws.Range("K98").Value = 651
$ws.Range("M98").Value = 847

$ws.Range("H107").Value = 766
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 766
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2298
$ws.Range("N107").Value = -6138
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 16380.333
$ws.Range("J86").Value = 16380.333
$ws.Range("L86").Value = 16380.333
$ws.Range("N86").Value = -18752.333

$ws.Range("H89").Value = 16380.333
$ws.Range("J89").Value = 16380.333
$ws.Range("L89").Value = 49140.999
$ws.Range("N89").Value = -60996.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1492.8823
$ws.Range("I7").Value = 1047.2727
$ws.Range("K7").Value = 1047.2727
$ws.Range("M7").Value = -935.2727

$ws.Range("H126").Value = 1492.8823
$ws.Range("I126").Value = 1047.2727
$ws.Range("K126").Value = 3141.8181
$ws.Range("M126").Value = -671.8181

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
